$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure text-like numeric values (e.g. "1.00", "408.04") are stored as text,
# matching the original inlineStr cell type, instead of being auto-converted
# to numbers (which would strip formatting like trailing zeros).
$textCells = @("D2", "D3", "D5", "D6", "D7", "D9", "D10", "D11", "D12", "D14", "D15", "D16", "D17", "D18", "D19", "D20", "D21", "D22", "D24", "D25", "D26", "D27", "D28", "D29", "D30", "D31", "D32", "D34", "D36", "D37", "D38", "D39", "D40", "D41", "D42", "D43", "D44", "D45", "D47", "D48", "D49", "D50", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated values from the latest crypto data refresh.
$ws.Range("D2").Value = '61.934.43'
$ws.Range("E2").Value = '  -0.04%  '
$ws.Range("D3").Value = '3.408.04'
$ws.Range("E3").Value = '  +0.08%  '
$ws.Range("E4").Value = '  +0.25%  '
$ws.Range("D5").Value = '408.04'
$ws.Range("E5").Value = '  +0.14%  '
$ws.Range("D6").Value = '128.38'
$ws.Range("E6").Value = '  -3.83%  '
$ws.Range("D7").Value = '0.633'
$ws.Range("E7").Value = '  +6.76%  '
$ws.Range("E8").Value = '  -0.06%  '
$ws.Range("D9").Value = '0.730'
$ws.Range("E9").Value = '  +7.10%  '
$ws.Range("D10").Value = '0.142'
$ws.Range("E10").Value = '  +17.57%  '
$ws.Range("D11").Value = '42.26'
$ws.Range("E11").Value = '  -0.89%  '
$ws.Range("D12").Value = '0.0000220'
$ws.Range("E12").Value = '  +69.66%  '
$ws.Range("E13").Value = '  -0.43%  '
$ws.Range("D14").Value = '3.960.00'
$ws.Range("E14").Value = '  +0.19%  '
$ws.Range("D15").Value = '8.89'
$ws.Range("E15").Value = '  +5.79%  '
$ws.Range("D16").Value = '20.69'
$ws.Range("E16").Value = '  +4.29%  '
$ws.Range("D17").Value = '3.407.26'
$ws.Range("E17").Value = '  -0.91%  '
$ws.Range("D18").Value = '12.09'
$ws.Range("E18").Value = '  +9.91%  '
$ws.Range("D19").Value = '1.06'
$ws.Range("E19").Value = '  +5.36%  '
$ws.Range("D20").Value = '61.897.08'
$ws.Range("E20").Value = '  -0.21%  '
$ws.Range("D21").Value = '410.49'
$ws.Range("E21").Value = '  +30.96%  '
$ws.Range("D22").Value = '89.17'
$ws.Range("E22").Value = '  +6.12%  '
$ws.Range("E23").Value = '  -0.75%  '
$ws.Range("D24").Value = '13.01'
$ws.Range("E24").Value = '  +1.32%  '
$ws.Range("D25").Value = '3.22'
$ws.Range("E25").Value = '  +1.79%  '
$ws.Range("D26").Value = '32.90'
$ws.Range("E26").Value = '  +11.33%  '
$ws.Range("D27").Value = '8.82'
$ws.Range("E27").Value = '  +8.27%  '
$ws.Range("D28").Value = '4.78'
$ws.Range("E28").Value = '  -0.19%  '
$ws.Range("D29").Value = '7.58'
$ws.Range("E29").Value = '  -0.16%  '
$ws.Range("D30").Value = '2.71'
$ws.Range("E30").Value = '  -3.06%  '
$ws.Range("B31").Value = 'Kaspa'
$ws.Range("C31").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D31").Value = '0.171'
$ws.Range("E31").Value = '  -1.94%  '
$ws.Range("B32").Value = 'Cosmos'
$ws.Range("C32").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D32").Value = '11.84'
$ws.Range("E32").Value = '  +4.02%  '
$ws.Range("E33").Value = '  -0.27%  '
$ws.Range("D34").Value = '42.53'
$ws.Range("E34").Value = '  -0.53%  '
$ws.Range("E35").Value = '  +0.79%  '
$ws.Range("D36").Value = '0.0496'
$ws.Range("E36").Value = '  +2.98%  '
$ws.Range("D37").Value = '53.99'
$ws.Range("E37").Value = '  +4.40%  '
$ws.Range("D38").Value = '1.00'
$ws.Range("E38").Value = '  -0.04%  '
$ws.Range("D39").Value = '3.34'
$ws.Range("E39").Value = '  -2.32%  '
$ws.Range("D40").Value = '0.133'
$ws.Range("E40").Value = '  +6.47%  '
$ws.Range("D41").Value = '2.90'
$ws.Range("E41").Value = '  -1.34%  '
$ws.Range("D42").Value = '0.310'
$ws.Range("E42").Value = '  +3.16%  '
$ws.Range("D43").Value = '141.62'
$ws.Range("D44").Value = '1.96'
$ws.Range("E44").Value = '  -1.23%  '
$ws.Range("D45").Value = '4.08'
$ws.Range("E45").Value = '  +1.31%  '
$ws.Range("E46").Value = '  +8.46%  '
$ws.Range("D47").Value = '16.58'
$ws.Range("E47").Value = '  -0.70%  '
$ws.Range("D48").Value = '21.74'
$ws.Range("E48").Value = '  +2.20%  '
$ws.Range("D49").Value = '2.106.33'
$ws.Range("E49").Value = '  -0.65%  '
$ws.Range("D50").Value = '2.38'
$ws.Range("E50").Value = '  +2.71%  '
$ws.Range("D51").Value = '0.131'
$ws.Range("E51").Value = '  +14.94%  '
